$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 654
$ws.Range("A653:L653").Copy($ws.Range("A654:L654"))
$ws.Range("M653").Copy($ws.Range("M654"))
$ws.Range("A654").Value = 45192.84447600694
$ws.Range("B654").Value = 'qwop7845@naver.com'
$ws.Range("C654").Value = '미래융합스쿨'
$ws.Range("D654").Value = 20236638
$ws.Range("E654").Value = '최시연'
$ws.Range("F654").Value = "'76:24"
$ws.Range("G654").Value = 0.2
$ws.Range("H654").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I654").Value = '38만 명'
$ws.Range("J654").Value = 0.151
$ws.Range("K654").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L654").Value = 'Red'
$ws.Range("M654").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 655
$ws.Range("A653:L653").Copy($ws.Range("A655:L655"))
$ws.Range("M653").Copy($ws.Range("M655"))
$ws.Range("A655").Value = 45192.87225399306
$ws.Range("B655").Value = 'h20221203@glab.ac.kr'
$ws.Range("C655").Value = '영어영문학과'
$ws.Range("D655").Value = 20221203
$ws.Range("E655").Value = '권민주'
$ws.Range("F655").Value = "'74:26"
$ws.Range("G655").Value = 0.2
$ws.Range("H655").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I655").Value = '952만 명'
$ws.Range("J655").Value = 0.059
$ws.Range("K655").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L655").Value = 'Red'
$ws.Range("M655").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

# Row 656
$ws.Range("A653:L653").Copy($ws.Range("A656:L656"))
$ws.Range("M653").Copy($ws.Range("M656"))
$ws.Range("A656").Value = 45192.88932383102
$ws.Range("B656").Value = 'audwlswlghd7@naver.com'
$ws.Range("C656").Value = '경영학과'
$ws.Range("D656").Value = 20192902
$ws.Range("E656").Value = '신명진'
$ws.Range("F656").Value = "'77:23"
$ws.Range("G656").Value = 0.2
$ws.Range("H656").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I656").Value = '779만 명'
$ws.Range("J656").Value = 0.059
$ws.Range("K656").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L656").Value = 'Red'
$ws.Range("M656").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

# Row 657
$ws.Range("A653:L653").Copy($ws.Range("A657:L657"))
$ws.Range("M653").Copy($ws.Range("M657"))
$ws.Range("A657").Value = 45192.89174630787
$ws.Range("B657").Value = 'rorita1191@naver.com'
$ws.Range("C657").Value = '간호학과'
$ws.Range("D657").Value = 20236225
$ws.Range("E657").Value = '김연슬'
$ws.Range("F657").Value = "'75:25"
$ws.Range("G657").Value = 0.25
$ws.Range("H657").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I657").Value = '166만 명'
$ws.Range("J657").Value = 0.374
$ws.Range("K657").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L657").Value = 'Red'
$ws.Range("M657").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 658
$ws.Range("A653:L653").Copy($ws.Range("A658:L658"))
$ws.Range("M653").Copy($ws.Range("M658"))
$ws.Range("A658").Value = 45192.901828807866
$ws.Range("B658").Value = 'sea36987412@gmail.com'
$ws.Range("C658").Value = '소프트웨어학부'
$ws.Range("D658").Value = 20235213
$ws.Range("E658").Value = '유세아'
$ws.Range("F658").Value = "'74:26"
$ws.Range("G658").Value = 0.2
$ws.Range("H658").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I658").Value = '952만 명'
$ws.Range("J658").Value = 0.059
$ws.Range("K658").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L658").Value = 'Red'
$ws.Range("M658").Value = '모름/무응답'

# Row 659
$ws.Range("A652:L652").Copy($ws.Range("A659:L659"))
$ws.Range("N652").Copy($ws.Range("N659"))
$ws.Range("A659").Value = 45192.91431248843
$ws.Range("B659").Value = 'yglee1357@naver.com'
$ws.Range("C659").Value = '생명과학과'
$ws.Range("D659").Value = 20193532
$ws.Range("E659").Value = '이윤구'
$ws.Range("F659").Value = "'74:26"
$ws.Range("G659").Value = 0.2
$ws.Range("H659").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I659").Value = '952만 명'
$ws.Range("J659").Value = 0.059
$ws.Range("K659").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L659").Value = 'Black'
$ws.Range("N659").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 660
$ws.Range("A652:L652").Copy($ws.Range("A660:L660"))
$ws.Range("N652").Copy($ws.Range("N660"))
$ws.Range("A660").Value = 45192.91745818287
$ws.Range("B660").Value = 'codmsrjf@naver.com'
$ws.Range("C660").Value = '콘텐츠IT'
$ws.Range("D660").Value = 20205253
$ws.Range("E660").Value = '정채은'
$ws.Range("F660").Value = "'77:23"
$ws.Range("G660").Value = 0.15
$ws.Range("H660").Value = '조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다'
$ws.Range("I660").Value = '779만 명'
$ws.Range("J660").Value = 0.151
$ws.Range("K660").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L660").Value = 'Black'
$ws.Range("N660").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 661
$ws.Range("A653:L653").Copy($ws.Range("A661:L661"))
$ws.Range("M653").Copy($ws.Range("M661"))
$ws.Range("A661").Value = 45192.91868228009
$ws.Range("B661").Value = 'laurano.first@gmail.com'
$ws.Range("C661").Value = '소프트웨어학부'
$ws.Range("D661").Value = 20235274
$ws.Range("E661").Value = '한서연'
$ws.Range("F661").Value = "'76:24"
$ws.Range("G661").Value = 0.25
$ws.Range("H661").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I661").Value = '166만 명'
$ws.Range("J661").Value = 0.151
$ws.Range("K661").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L661").Value = 'Red'
$ws.Range("M661").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 662
$ws.Range("A653:L653").Copy($ws.Range("A662:L662"))
$ws.Range("M653").Copy($ws.Range("M662"))
$ws.Range("A662").Value = 45192.92362686343
$ws.Range("B662").Value = 'wlsgml9808@naver.com'
$ws.Range("C662").Value = '환경생명공학과'
$ws.Range("D662").Value = 20173728
$ws.Range("E662").Value = '이진희'
$ws.Range("F662").Value = "'74:26"
$ws.Range("G662").Value = 0.1
$ws.Range("H662").Value = '조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다'
$ws.Range("I662").Value = '952만 명'
$ws.Range("J662").Value = 0.002
$ws.Range("K662").Value = '중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다'
$ws.Range("L662").Value = 'Red'
$ws.Range("M662").Value = '모름/무응답'

# Row 663
$ws.Range("A653:L653").Copy($ws.Range("A663:L663"))
$ws.Range("M653").Copy($ws.Range("M663"))
$ws.Range("A663").Value = 45192.94106090278
$ws.Range("B663").Value = 'snowy0601@naver.com'
$ws.Range("C663").Value = '미디어스쿨'
$ws.Range("D663").Value = 20232582
$ws.Range("E663").Value = '차연우'
$ws.Range("F663").Value = "'76:24"
$ws.Range("G663").Value = 0.2
$ws.Range("H663").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I663").Value = '952만 명'
$ws.Range("J663").Value = 0.059
$ws.Range("K663").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L663").Value = 'Red'
$ws.Range("M663").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 664
$ws.Range("A653:L653").Copy($ws.Range("A664:L664"))
$ws.Range("M653").Copy($ws.Range("M664"))
$ws.Range("A664").Value = 45192.95097297453
$ws.Range("B664").Value = '20232593@hallym.ac.kr'
$ws.Range("C664").Value = '미디어스쿨'
$ws.Range("D664").Value = 20232593
$ws.Range("E664").Value = '김나영'
$ws.Range("F664").Value = "'78:22"
$ws.Range("G664").Value = 0.15
$ws.Range("H664").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I664").Value = '952만 명'
$ws.Range("J664").Value = 0.151
$ws.Range("K664").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L664").Value = 'Red'
$ws.Range("M664").Value = '모름/무응답'

# Row 665
$ws.Range("A652:L652").Copy($ws.Range("A665:L665"))
$ws.Range("N652").Copy($ws.Range("N665"))
$ws.Range("A665").Value = 45192.956725613425
$ws.Range("B665").Value = 'a01093819049@gmail.com'
$ws.Range("C665").Value = '광고홍보학과'
$ws.Range("D665").Value = 20212608
$ws.Range("E665").Value = '김서영'
$ws.Range("F665").Value = "'75:25"
$ws.Range("G665").Value = 0.2
$ws.Range("H665").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I665").Value = '166만 명'
$ws.Range("J665").Value = 0.374
$ws.Range("K665").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L665").Value = 'Black'
$ws.Range("N665").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

# Row 666
$ws.Range("A653:L653").Copy($ws.Range("A666:L666"))
$ws.Range("M653").Copy($ws.Range("M666"))
$ws.Range("A666").Value = 45192.96557398148
$ws.Range("B666").Value = 'wusl0327@naver.com'
$ws.Range("C666").Value = '청각학전공'
$ws.Range("D666").Value = 20193934
$ws.Range("E666").Value = '신지연'
$ws.Range("F666").Value = "'75:25"
$ws.Range("G666").Value = 0.2
$ws.Range("H666").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I666").Value = '952만 명'
$ws.Range("J666").Value = 0.002
$ws.Range("K666").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L666").Value = 'Red'
$ws.Range("M666").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 667
$ws.Range("A652:L652").Copy($ws.Range("A667:L667"))
$ws.Range("N652").Copy($ws.Range("N667"))
$ws.Range("A667").Value = 45192.97250038195
$ws.Range("B667").Value = 'sungyeon0803@gmail.com'
$ws.Range("C667").Value = '정치행정학과'
$ws.Range("D667").Value = 20232437
$ws.Range("E667").Value = '최성연'
$ws.Range("F667").Value = "'74:26"
$ws.Range("G667").Value = 0.2
$ws.Range("H667").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I667").Value = '952만 명'
$ws.Range("J667").Value = 0.059
$ws.Range("K667").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("L667").Value = 'Black'
$ws.Range("N667").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 668
$ws.Range("A653:L653").Copy($ws.Range("A668:L668"))
$ws.Range("M653").Copy($ws.Range("M668"))
$ws.Range("A668").Value = 45192.97574122685
$ws.Range("B668").Value = 'jayean0715@naver.com'
$ws.Range("C668").Value = '언어청각학부'
$ws.Range("D668").Value = 20233918
$ws.Range("E668").Value = '김자연'
$ws.Range("F668").Value = "'74:26"
$ws.Range("G668").Value = 0.2
$ws.Range("H668").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I668").Value = '952만 명'
$ws.Range("J668").Value = 0.059
$ws.Range("K668").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L668").Value = 'Red'
$ws.Range("M668").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 669
$ws.Range("A653:L653").Copy($ws.Range("A669:L669"))
$ws.Range("M653").Copy($ws.Range("M669"))
$ws.Range("A669").Value = 45192.98677287037
$ws.Range("B669").Value = 'sulnak159@gmail.com'
$ws.Range("C669").Value = '일본학과'
$ws.Range("D669").Value = 20181605
$ws.Range("E669").Value = '김남준'
$ws.Range("F669").Value = "'74:26"
$ws.Range("G669").Value = 0.2
$ws.Range("H669").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I669").Value = '952만 명'
$ws.Range("J669").Value = 0.059
$ws.Range("K669").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L669").Value = 'Red'
$ws.Range("M669").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 670
$ws.Range("A652:L652").Copy($ws.Range("A670:L670"))
$ws.Range("N652").Copy($ws.Range("N670"))
$ws.Range("A670").Value = 45192.98951328704
$ws.Range("B670").Value = 'richhjin@naver.com'
$ws.Range("C670").Value = '환경생명공학과'
$ws.Range("D670").Value = 20233732
$ws.Range("E670").Value = '이현진'
$ws.Range("F670").Value = "'74:26"
$ws.Range("G670").Value = 0.2
$ws.Range("H670").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I670").Value = '952만 명'
$ws.Range("J670").Value = 0.059
$ws.Range("K670").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L670").Value = 'Black'
$ws.Range("N670").Value = '모름/무응답'

# Row 671
$ws.Range("A652:L652").Copy($ws.Range("A671:L671"))
$ws.Range("N652").Copy($ws.Range("N671"))
$ws.Range("A671").Value = 45192.99965989584
$ws.Range("B671").Value = 'sjh8358000@naver.com'
$ws.Range("C671").Value = '언어청각학부'
$ws.Range("D671").Value = 20233934
$ws.Range("E671").Value = '손지후'
$ws.Range("F671").Value = "'78:22"
$ws.Range("G671").Value = 0.2
$ws.Range("H671").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I671").Value = '38만 명'
$ws.Range("J671").Value = 0.151
$ws.Range("K671").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L671").Value = 'Black'
$ws.Range("N671").Value = '모름/무응답'

# Row 672
$ws.Range("A652:L652").Copy($ws.Range("A672:L672"))
$ws.Range("N652").Copy($ws.Range("N672"))
$ws.Range("A672").Value = 45193.00692424769
$ws.Range("B672").Value = 'was193567@gmail.com'
$ws.Range("C672").Value = '반도체디스플레이스쿨'
$ws.Range("D672").Value = 20233303
$ws.Range("E672").Value = '권유진'
$ws.Range("F672").Value = "'75:25"
$ws.Range("G672").Value = 0.2
$ws.Range("H672").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I672").Value = '952만 명'
$ws.Range("J672").Value = 0.059
$ws.Range("K672").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("L672").Value = 'Black'
$ws.Range("N672").Value = '모름/무응답'

# Row 673
$ws.Range("A652:L652").Copy($ws.Range("A673:L673"))
$ws.Range("N652").Copy($ws.Range("N673"))
$ws.Range("A673").Value = 45193.013319687496
$ws.Range("B673").Value = 'leah0820@naver.com'
$ws.Range("C673").Value = '경영학과'
$ws.Range("D673").Value = 20232981
$ws.Range("E673").Value = '심채린'
$ws.Range("F673").Value = "'77:23"
$ws.Range("G673").Value = 0.15
$ws.Range("H673").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I673").Value = '166만 명'
$ws.Range("J673").Value = 0.151
$ws.Range("K673").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L673").Value = 'Black'
$ws.Range("N673").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 674
$ws.Range("A653:L653").Copy($ws.Range("A674:L674"))
$ws.Range("M653").Copy($ws.Range("M674"))
$ws.Range("A674").Value = 45193.0292554051
$ws.Range("B674").Value = 'skysun0510@naver.com'
$ws.Range("C674").Value = '간호학과'
$ws.Range("D674").Value = 20236268
$ws.Range("E674").Value = '유중선'
$ws.Range("F674").Value = "'77:23"
$ws.Range("G674").Value = 0.2
$ws.Range("H674").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I674").Value = '166만 명'
$ws.Range("J674").Value = 0.151
$ws.Range("K674").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L674").Value = 'Red'
$ws.Range("M674").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 675
$ws.Range("A653:L653").Copy($ws.Range("A675:L675"))
$ws.Range("M653").Copy($ws.Range("M675"))
$ws.Range("A675").Value = 45193.035610011575
$ws.Range("B675").Value = 'daysyad@naver.com'
$ws.Range("C675").Value = '미디어스쿨'
$ws.Range("D675").Value = 20232586
$ws.Range("E675").Value = '최소연'
$ws.Range("F675").Value = "'74:26"
$ws.Range("G675").Value = 0.2
$ws.Range("H675").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I675").Value = '952만 명'
$ws.Range("J675").Value = 0.059
$ws.Range("K675").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L675").Value = 'Red'
$ws.Range("M675").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 676
$ws.Range("A653:L653").Copy($ws.Range("A676:L676"))
$ws.Range("M653").Copy($ws.Range("M676"))
$ws.Range("A676").Value = 45193.05239094907
$ws.Range("B676").Value = 'jeongyeon256@gmail.com'
$ws.Range("C676").Value = '사회학과'
$ws.Range("D676").Value = 20232232
$ws.Range("E676").Value = '이정연'
$ws.Range("F676").Value = "'76:24"
$ws.Range("G676").Value = 0.15
$ws.Range("H676").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I676").Value = '166만 명'
$ws.Range("J676").Value = 0.151
$ws.Range("K676").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L676").Value = 'Red'
$ws.Range("M676").Value = '모름/무응답'

# Row 677
$ws.Range("A653:L653").Copy($ws.Range("A677:L677"))
$ws.Range("M653").Copy($ws.Range("M677"))
$ws.Range("A677").Value = 45193.06368795139
$ws.Range("B677").Value = 'jangyc3972@gmail.com'
$ws.Range("C677").Value = '인공지능융합학부'
$ws.Range("D677").Value = 20236771
$ws.Range("E677").Value = '장윤채'
$ws.Range("F677").Value = "'76:24"
$ws.Range("G677").Value = 0.15
$ws.Range("H677").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I677").Value = '166만 명'
$ws.Range("J677").Value = 0.059
$ws.Range("K677").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("L677").Value = 'Red'
$ws.Range("M677").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

# Row 678
$ws.Range("A652:L652").Copy($ws.Range("A678:L678"))
$ws.Range("N652").Copy($ws.Range("N678"))
$ws.Range("A678").Value = 45193.09273060186
$ws.Range("B678").Value = 'hankyo777@naver.com'
$ws.Range("C678").Value = '데이터사이언스학부'
$ws.Range("D678").Value = 20233244
$ws.Range("E678").Value = '장성주'
$ws.Range("F678").Value = "'74:26"
$ws.Range("G678").Value = 0.2
$ws.Range("H678").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I678").Value = '779만 명'
$ws.Range("J678").Value = 0.151
$ws.Range("K678").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L678").Value = 'Black'
$ws.Range("N678").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 679
$ws.Range("A653:L653").Copy($ws.Range("A679:L679"))
$ws.Range("M653").Copy($ws.Range("M679"))
$ws.Range("A679").Value = 45193.132103020835
$ws.Range("B679").Value = 'juliaqkrwpdm@hallym.ac.kr'
$ws.Range("C679").Value = '간호학과'
$ws.Range("D679").Value = 20236241
$ws.Range("E679").Value = '박제은'
$ws.Range("F679").Value = "'74:26"
$ws.Range("G679").Value = 0.2
$ws.Range("H679").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I679").Value = '952만 명'
$ws.Range("J679").Value = 0.059
$ws.Range("K679").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L679").Value = 'Red'
$ws.Range("M679").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 680
$ws.Range("A653:L653").Copy($ws.Range("A680:L680"))
$ws.Range("M653").Copy($ws.Range("M680"))
$ws.Range("A680").Value = 45193.1338143287
$ws.Range("B680").Value = 'changyw6729@naver.com'
$ws.Range("C680").Value = '미디어스쿨'
$ws.Range("D680").Value = 20232573
$ws.Range("E680").Value = '장연욱'
$ws.Range("F680").Value = "'74:26"
$ws.Range("G680").Value = 0.2
$ws.Range("H680").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I680").Value = '952만 명'
$ws.Range("J680").Value = 0.059
$ws.Range("K680").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L680").Value = 'Red'
$ws.Range("M680").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 681
$ws.Range("A652:L652").Copy($ws.Range("A681:L681"))
$ws.Range("N652").Copy($ws.Range("N681"))
$ws.Range("A681").Value = 45193.14376489584
$ws.Range("B681").Value = 'choijoonhyuk1412@naver.com'
$ws.Range("C681").Value = '미디어스쿨'
$ws.Range("D681").Value = 20232588
$ws.Range("E681").Value = '최준혁'
$ws.Range("F681").Value = "'74:26"
$ws.Range("G681").Value = 0.2
$ws.Range("H681").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I681").Value = '952만 명'
$ws.Range("J681").Value = 0.059
$ws.Range("K681").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L681").Value = 'Black'
$ws.Range("N681").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 682
$ws.Range("A652:L652").Copy($ws.Range("A682:L682"))
$ws.Range("N652").Copy($ws.Range("N682"))
$ws.Range("A682").Value = 45193.158318518515
$ws.Range("B682").Value = 'leegijae040209@gmail.com'
$ws.Range("C682").Value = '인공지능융합학부'
$ws.Range("D682").Value = 20236761
$ws.Range("E682").Value = '이윤'
$ws.Range("F682").Value = "'74:26"
$ws.Range("G682").Value = 0.2
$ws.Range("H682").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I682").Value = '952만 명'
$ws.Range("J682").Value = 0.059
$ws.Range("K682").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L682").Value = 'Black'
$ws.Range("N682").Value = '모름/무응답'

# Row 683
$ws.Range("A653:L653").Copy($ws.Range("A683:L683"))
$ws.Range("M653").Copy($ws.Range("M683"))
$ws.Range("A683").Value = 45193.178574606485
$ws.Range("B683").Value = 'gilh3408@gmail.com'
$ws.Range("C683").Value = '데이터사이언스학부'
$ws.Range("D683").Value = 20233204
$ws.Range("E683").Value = '길혜균'
$ws.Range("F683").Value = "'74:26"
$ws.Range("G683").Value = 0.2
$ws.Range("H683").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I683").Value = '952만 명'
$ws.Range("J683").Value = 0.059
$ws.Range("K683").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L683").Value = 'Red'
$ws.Range("M683").Value = '모름/무응답'

# Row 684
$ws.Range("A653:L653").Copy($ws.Range("A684:L684"))
$ws.Range("M653").Copy($ws.Range("M684"))
$ws.Range("A684").Value = 45193.24284328704
$ws.Range("B684").Value = 'limyoon0725@daum.net'
$ws.Range("C684").Value = '사회복지학과'
$ws.Range("D684").Value = 20222347
$ws.Range("E684").Value = '임윤서'
$ws.Range("F684").Value = "'76:24"
$ws.Range("G684").Value = 0.2
$ws.Range("H684").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I684").Value = '779만 명'
$ws.Range("J684").Value = 0.151
$ws.Range("K684").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L684").Value = 'Red'
$ws.Range("M684").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 685
$ws.Range("A653:L653").Copy($ws.Range("A685:L685"))
$ws.Range("M653").Copy($ws.Range("M685"))
$ws.Range("A685").Value = 45193.28391138889
$ws.Range("B685").Value = 'dennis121691@gmail.com'
$ws.Range("C685").Value = '  데이터테크전공'
$ws.Range("D685").Value = 20203224
$ws.Range("E685").Value = '서동욱'
$ws.Range("F685").Value = "'74:26"
$ws.Range("G685").Value = 0.2
$ws.Range("H685").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I685").Value = '952만 명'
$ws.Range("J685").Value = 0.059
$ws.Range("K685").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L685").Value = 'Red'
$ws.Range("M685").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

